$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196, shifting existing rows 196:270 down to 197:271.
$ws.Rows("196:196").Insert()

# Populate the newly inserted row 196 with the new weekly record's data.
$ws.Cells.Item(196, 1).Value = 4
$ws.Cells.Item(196, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(196, 3).Value = 'Los Lagos'
$ws.Cells.Item(196, 4).Value = 44468
$ws.Cells.Item(196, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(196, 5).Value = 10
$ws.Cells.Item(196, 6).Value = 'Fruta'
$ws.Cells.Item(196, 7).Value = 100102
$ws.Cells.Item(196, 8).Value = 'Cítricos'
$ws.Cells.Item(196, 9).Value = 100102003
$ws.Cells.Item(196, 10).Value = 'Limón'
$ws.Cells.Item(196, 11).Value = 'Sin especificar'
$ws.Cells.Item(196, 12).Value = '1a amarillo'
$ws.Cells.Item(196, 13).Value = 300
$ws.Cells.Item(196, 14).Value = 9000
$ws.Cells.Item(196, 15).Value = 9000
$ws.Cells.Item(196, 16).Value = 9000
$ws.Cells.Item(196, 17).Value = '$/malla 18 kilos'
$ws.Cells.Item(196, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(196, 19).Value = 500
$ws.Cells.Item(196, 20).Value = 18
